$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    31  = @("sd", "Statement-non-opinion")
    34  = @("sd", "Statement-non-opinion")
    35  = @("%", "Uninterpretable")
    38  = @("sd", "Statement-non-opinion")
    48  = @("b", "Acknowledge (Backchannel)")
    52  = @("%", "Uninterpretable")
    60  = @("aa", "Agree/Accept")
    66  = @("sd", "Statement-non-opinion")
    75  = @("aa", "Agree/Accept")
    79  = @("aa", "Agree/Accept")
    83  = @("aa", "Agree/Accept")
    92  = @("sd", "Statement-non-opinion")
    144 = @("%", "Uninterpretable")
    180 = @("sv", "Statement-opinion")
    187 = @("sd", "Statement-non-opinion")
    189 = @("aa", "Agree/Accept")
    198 = @("%", "Uninterpretable")
    202 = @("%", "Uninterpretable")
    233 = @("aa", "Agree/Accept")
    283 = @("sd", "Statement-non-opinion")
    292 = @("sd", "Statement-non-opinion")
    295 = @("sd", "Statement-non-opinion")
    298 = @("sd", "Statement-non-opinion")
    299 = @("b", "Acknowledge (Backchannel)")
    306 = @("sv", "Statement-opinion")
    324 = @("b", "Acknowledge (Backchannel)")
    325 = @("sd", "Statement-non-opinion")
    365 = @("b", "Acknowledge (Backchannel)")
    376 = @("sd", "Statement-non-opinion")
    379 = @("sd", "Statement-non-opinion")
    381 = @("b", "Acknowledge (Backchannel)")
    403 = @("sv", "Statement-opinion")
    405 = @("sd", "Statement-non-opinion")
    415 = @("b", "Acknowledge (Backchannel)")
    418 = @("sv", "Statement-opinion")
    420 = @("sd", "Statement-non-opinion")
    428 = @("sd", "Statement-non-opinion")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
